$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2046979865771812
$ws.Range("C2").Value = 0.5335570469798657
$ws.Range("J2").Value = 0.01677852348993289
$ws.Range("P2").Value = 0.1409395973154362
$ws.Range("S2").Value = 0.1040268456375839
$ws.Range("B3").Value = 0.006211180124223602
$ws.Range("C3").Value = 0.0124223602484472
$ws.Range("J3").Value = 0.02484472049689441
$ws.Range("P3").Value = 0.8136645962732919
$ws.Range("S3").Value = 0.1428571428571428
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1875
$ws.Range("B6").Value = 0.05076142131979695
$ws.Range("D6").Value = 0.01015228426395939
$ws.Range("F6").Value = 0.05076142131979695
$ws.Range("J6").Value = 0.233502538071066
$ws.Range("O6").Value = 0.02030456852791878
$ws.Range("Q6").Value = 0.2436548223350254
$ws.Range("R6").Value = 0.07106598984771574
$ws.Range("S6").Value = 0.3197969543147208
$ws.Range("B7").Value = 0.0718954248366013
$ws.Range("D7").Value = 0.03267973856209151
$ws.Range("F7").Value = 0.0392156862745098
$ws.Range("J7").Value = 0.130718954248366
$ws.Range("O7").Value = 0.03267973856209151
$ws.Range("Q7").Value = 0.2091503267973856
$ws.Range("R7").Value = 0.08496732026143791
$ws.Range("S7").Value = 0.3986928104575164
$ws.Range("B8").Value = 0.102803738317757
$ws.Range("D8").Value = 0.02803738317757009
$ws.Range("F8").Value = 0.06775700934579439
$ws.Range("J8").Value = 0.1051401869158878
$ws.Range("O8").Value = 0.01401869158878505
$ws.Range("Q8").Value = 0.2593457943925234
$ws.Range("R8").Value = 0.09345794392523364
$ws.Range("S8").Value = 0.3294392523364486
$ws.Range("B9").Value = 0.1219512195121951
$ws.Range("D9").Value = 0.01951219512195122
$ws.Range("F9").Value = 0.06829268292682927
$ws.Range("J9").Value = 0.1317073170731707
$ws.Range("O9").Value = 0.02439024390243903
$ws.Range("Q9").Value = 0.2341463414634146
$ws.Range("R9").Value = 0.07317073170731707
$ws.Range("S9").Value = 0.3268292682926829
$ws.Range("B10").Value = 0.1082474226804124
$ws.Range("D10").Value = 0.02945508100147275
$ws.Range("F10").Value = 0.06553755522827688
$ws.Range("J10").Value = 0.1207658321060383
$ws.Range("O10").Value = 0.02135493372606775
$ws.Range("Q10").Value = 0.2842415316642121
$ws.Range("R10").Value = 0.07658321060382917
$ws.Range("S10").Value = 0.2938144329896907
$ws.Range("G11").Value = 0.1410788381742739
$ws.Range("J11").Value = 0.07468879668049792
$ws.Range("K11").Value = 0.1950207468879668
$ws.Range("L11").Value = 0.5767634854771784
$ws.Range("S11").Value = 0.01244813278008299
$ws.Range("G12").Value = 0.7266187050359713
$ws.Range("J12").Value = 0.2302158273381295
$ws.Range("K12").Value = 0.02877697841726619
$ws.Range("S12").Value = 0.01438848920863309
$ws.Range("F15").Value = 0.004166666666666667
$ws.Range("H15").Value = 0.1208333333333333
$ws.Range("I15").Value = 0.05833333333333333
$ws.Range("J15").Value = 0.4208333333333333
$ws.Range("K15").Value = 0.03333333333333333
$ws.Range("M15").Value = 0.008333333333333333
$ws.Range("O15").Value = 0.075
$ws.Range("S15").Value = 0.2791666666666667
$ws.Range("F16").Value = 0.02304147465437788
$ws.Range("H16").Value = 0.1889400921658986
$ws.Range("I16").Value = 0.06912442396313365
$ws.Range("J16").Value = 0.4331797235023042
$ws.Range("K16").Value = 0.07834101382488479
$ws.Range("M16").Value = 0.01382488479262673
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.1290322580645161
$ws.Range("F17").Value = 0.009646302250803859
$ws.Range("H17").Value = 0.1446945337620579
$ws.Range("I17").Value = 0.09646302250803858
$ws.Range("J17").Value = 0.4678456591639871
$ws.Range("K17").Value = 0.09485530546623794
$ws.Range("M17").Value = 0.01446945337620579
$ws.Range("N17").Value = 0.001607717041800643
$ws.Range("O17").Value = 0.05787781350482315
$ws.Range("S17").Value = 0.112540192926045
$ws.Range("F18").Value = 0.01063829787234043
$ws.Range("H18").Value = 0.2127659574468085
$ws.Range("I18").Value = 0.0851063829787234
$ws.Range("J18").Value = 0.4361702127659575
$ws.Range("K18").Value = 0.06382978723404255
$ws.Range("M18").Value = 0.005319148936170213
$ws.Range("O18").Value = 0.07446808510638298
$ws.Range("S18").Value = 0.1117021276595745
$ws.Range("F19").Value = 0.007434944237918215
$ws.Range("H19").Value = 0.20817843866171
$ws.Range("I19").Value = 0.09479553903345725
$ws.Range("J19").Value = 0.4024163568773234
$ws.Range("K19").Value = 0.08828996282527882
$ws.Range("M19").Value = 0.01672862453531598
$ws.Range("N19").Value = 0.0009293680297397769
$ws.Range("O19").Value = 0.08085501858736059
$ws.Range("S19").Value = 0.1003717472118959
